$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "56×88="
$t.Cell(1,2).Range.Text = "35×13="
$t.Cell(1,3).Range.Text = "16×28="
$t.Cell(1,4).Range.Text = "54×92="
$t.Cell(1,5).Range.Text = "89×16="
$t.Cell(5,1).Range.Text = "60×98="
$t.Cell(5,2).Range.Text = "21×49="
$t.Cell(5,3).Range.Text = "99×65="
$t.Cell(5,4).Range.Text = "27×40="
$t.Cell(5,5).Range.Text = "99×19="
$t.Cell(10,1).Range.Text = "17×79="
$t.Cell(10,2).Range.Text = "89×11="
$t.Cell(10,3).Range.Text = "28×58="
$t.Cell(10,4).Range.Text = "30×57="
$t.Cell(10,5).Range.Text = "70×28="
$t.Cell(15,1).Range.Text = "17×28="
$t.Cell(15,2).Range.Text = "94×39="
$t.Cell(15,3).Range.Text = "50×38="
$t.Cell(15,4).Range.Text = "51×64="
$t.Cell(15,5).Range.Text = "46×69="
$t.Cell(20,1).Range.Text = "68×92="
$t.Cell(20,2).Range.Text = "29×21="
$t.Cell(20,3).Range.Text = "98×26="
$t.Cell(20,4).Range.Text = "23×78="
$t.Cell(20,5).Range.Text = "85×77="


Write-Host "Done updating table cells."
